# Replace all occurrences of "OIE" with "WOAH" across the relevant text
# cells of the workbook (Sheet 1 content column and References column),
# matching the upstream rename of the organisation acronym.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$sheet2 = $wb.Worksheets.Item("References")

# Cells on "Sheet 1" (column E = Content) that mention OIE.
$sheet1Cells = @("E5", "E6", "E7", "E14", "E17", "E31", "E44", "E53", "E65", "E77", "E158")

foreach ($addr in $sheet1Cells) {
    $cell = $sheet1.Range($addr)
    $text = $cell.Text
    $cell.Value = $text.Replace("OIE", "WOAH")
}

# Cells on "References" (column C = Paper) that mention OIE.
$sheet2Cells = @("C2", "C5", "C8", "C10")

foreach ($addr in $sheet2Cells) {
    $cell = $sheet2.Range($addr)
    $text = $cell.Text
    $cell.Value = $text.Replace("OIE", "WOAH")
}
